$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column I (RF) for rows 20 through 44 with the new RF value
$newValue = 5.292199999999999
for ($row = 20; $row -le 44; $row++) {
    $ws.Range("I$row").Value = $newValue
}
